# "added keyword tryScroll function"
# Insert a new TestSteps row for the "tryScroll" keyword right before the
# existing TS_41 ("slide") step, renumbering the TS_NN labels of the two
# steps that shift down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestSteps")

# Insert a blank row at 42: old row 42 (TS_41 / slide) becomes row 43,
# old row 43 (TS_42 / verify text) becomes row 44.
$ws.Rows.Item(42).Insert()

# New row 42: the tryScroll keyword step (reuses the TS_41 step id).
$ws.Cells.Item(42, 1).Value = "ServiceAdv"
$ws.Cells.Item(42, 2).Value = "TS_41"
$ws.Cells.Item(42, 3).Value = "scroll down"
$ws.Cells.Item(42, 5).Value = "tryScroll"

# The old TS_42 step (now row 44) becomes TS_43.
$ws.Cells.Item(44, 2).Value = "TS_43"

$ws.Cells.Item(42, 4).Value = "//*[@id=""sAdvantageOwl""]/div[1]/div/div[1]/div/div[2]/div[1]/table/tbody/tr[1]/td[2]/b"
$ws.Cells.Item(42, 7).Value = "PASS"

# The old TS_41 step (now row 43) becomes TS_42.
$ws.Cells.Item(43, 2).Value = "TS_42"

# Row 42 has no TestData value for this step; reset the inherited style so
# the cell serializes with no number format / style override.
$ws.Cells.Item(42, 6).Style = "Normal"

# Leave the selection where the author left it after the edit.
$null = $ws.Range("D42").Select()
